# Insert a new first data column (D) on the CMLS sheet for the new quarter
# (period ending 2018-09-30), shifting the previously-existing quarters
# (old D:K, ending 2018-06-30 .. 2016-09-30) one column to the right (E:L).
# Then populate the new column D with the newly reported quarter's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CMLS")

# 1) Insert a new column before D; existing D:K data shifts to E:L.
$ws.Columns("D").Insert()

# 2) The freshly inserted column starts out with default/general formatting.
#    Clone the number formats from column E (the old column D, now shifted)
#    onto the new column D so dates/numbers render the same way.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Fill in the new quarter's values (period ending 2018-09-30).
#    Dates are set as raw serial numbers (43373 = 2018-09-30) so the
#    pasted-in custom date format isn't clobbered by date auto-detection.
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 282300
$ws.Range("D9").Value = 98500
$ws.Range("D10").Value = 183800
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 14100
$ws.Range("D17").Value = 238900
$ws.Range("D18").Value = 43400
$ws.Range("D20").Value = -3500
$ws.Range("D21").Value = 54100
$ws.Range("D22").Value = 22100
$ws.Range("D23").Value = 17800
$ws.Range("D24").Value = 5100
$ws.Range("D26").Value = 12700
$ws.Range("D27").Value = 12700
$ws.Range("D29").Value = "NA"
$ws.Range("D32").Value = 3500
$ws.Range("D33").Value = 12700
$ws.Range("D35").Value = 12700

$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 54000
$ws.Range("D43").Value = 240500
$ws.Range("D45").Value = 125500
$ws.Range("D46").Value = 420000
$ws.Range("D48").Value = 235700
$ws.Range("D49").Value = 1137400
$ws.Range("D52").Value = 16200
$ws.Range("D54").Value = 1809400
$ws.Range("D57").Value = 96900
$ws.Range("D58").Value = 13000
$ws.Range("D60").Value = 109900
$ws.Range("D61").Value = 1283800
$ws.Range("D62").Value = 71300
$ws.Range("D66").Value = 1464900
$ws.Range("D72").Value = 17700
$ws.Range("D76").Value = 344500

$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 12700
$ws.Range("D83").Value = 14100
$ws.Range("D89").Value = 4200
$ws.Range("D91").Value = -5900
$ws.Range("D94").Value = -5900
$ws.Range("D100").Value = -3300
$ws.Range("D102").Value = -5000
